$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.176.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.07%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.802.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.24%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.38%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "163.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.10%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.799.18"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.28%  "

# Row 8
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$ws.Range("E9").Value = "  -0.40%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.170"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.95%  "

# Row 11
$ws.Range("E11").Value = "  -1.09%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.49%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.28"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.96%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000246"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.04%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.429.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.13%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.790.86"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.75%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.246.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.21%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.48%  "

# Row 19
$ws.Range("E19").Value = "  -0.18%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.13%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.32%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "488.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.55%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.723"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.90%  "

# Row 24
$ws.Range("E24").Value = "  -2.11%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.91%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.39%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.25"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.05%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.44%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.08%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.97"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.64%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.07"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.54%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.94%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.943.53"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.10%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.78"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.45%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.741.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.44%  "

# Row 36
$ws.Range("E36").Value = "  -1.98%  "

# Row 37
$ws.Range("E37").Value = "  +5.62%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.07%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.93"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.46%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.03%  "

# Row 41
$ws.Range("E41").Value = "  -0.59%  "

# Row 42
$ws.Range("E42").Value = "  -0.37%  "

# Row 43
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.21%  "

# Row 44
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.62"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.02%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "422.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.53%  "

# Row 46
$ws.Range("E46").Value = "  +0.00%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.23%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.825.81"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.45%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.15%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "39.81"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.61%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0351"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.29%  "
